$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.929.70'
$ws.Range("E2").Value = '  +5.93%  '

$ws.Range("D3").Value = '2.331.69'
$ws.Range("E3").Value = '  +4.92%  '

$ws.Range("E4").Value = '  -0.55%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.97%  '

$ws.Range("E7").Value = '  +4.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.541'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.23%  '

$ws.Range("E13").Value = '  +0.89%  '

$ws.Range("D14").Value = '2.686.02'
$ws.Range("E14").Value = '  +4.95%  '

$ws.Range("D15").Value = '2.328.94'
$ws.Range("E15").Value = '  +6.01%  '

$ws.Range("E16").Value = '  +5.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.58%  '

$ws.Range("D18").Value = '46.814.61'
$ws.Range("E18").Value = '  +6.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +21.59%  '

$ws.Range("D20").Value = '0.0₃0958'
$ws.Range("E20").Value = '  +5.94%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.80%  '

$ws.Range("E24").Value = '  +4.16%  '

$ws.Range("E25").Value = '  +5.25%  '

$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '42.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +15.92%  '

$ws.Range("E28").Value = '  +2.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0815'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '147.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.07%  '

$ws.Range("E34").Value = '  +2.37%  '

$ws.Range("E35").Value = '  +6.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.114'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.119'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.53%  '

$ws.Range("E38").Value = '  +6.13%  '

$ws.Range("E39").Value = '  +12.22%  '

$ws.Range("E40").Value = '  +8.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.95%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +20.22%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.72%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.814.65'
$ws.Range("E45").Value = '  +3.86%  '

$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +18.33%  '

$ws.Range("E47").Value = '  +8.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '74.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '98.95'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.36%  '
